$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / update existing rows and add a new row for the "Warhorse Blend (12oz)" product
$ws.Range("A1").Value = "Copper Horse - Rumble Pony (12oz)"
$ws.Range("B1").Value = 6

$ws.Range("A2").Value = "Copper Horse - Clocktower Espresso (12oz)"
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "Copper Horse - Carriage House Blend (12oz)"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "Copper Horse - Warhorse Blend (12oz)"
$ws.Range("B4").Value = 6

$ws.Range("A5").Value = "Copper Horse - Warhorse 5lb"
$ws.Range("B5").Value = 8
